# Update Name of Algo
# Applies the numeric value updates described in the commit diff to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value  = 6.328699999999999
$ws.Range("A3").Value  = -21.46130000000002
$ws.Range("B5").Value  = 4.618400000000005
$ws.Range("C5").Value  = -13.574
$ws.Range("D7").Value  = -7.081199999999996
$ws.Range("C9").Value  = -11.85960000000001
$ws.Range("C11").Value = -13.4644
$ws.Range("D11").Value = -8.022399999999998
$ws.Range("A14").Value = -20.51259999999997
$ws.Range("A16").Value = -21.48990000000002
$ws.Range("B16").Value = 5.940499999999995
$ws.Range("C17").Value = -11.3567
$ws.Range("D19").Value = -8.80229999999999
$ws.Range("A21").Value = -21.26300000000001
$ws.Range("C21").Value = -11.06109999999999
$ws.Range("D21").Value = -7.865800000000007
$ws.Range("A23").Value = -21.52350000000002
$ws.Range("A25").Value = -22.44400000000003
